$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. "As we work harder" -> "As we work a little harder"
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("s we work harder", $false, $false, $false, $false, $false, $true, 1, $false, "s we work a little harder", 2) | Out-Null

# ---------------------------------------------------------------------------
# 2. LT1 description rewrite:
#    "a point where we start to accumulate lactate in the legs, this point is
#     known as"
#    -> "a point where lactate levels increase above the baseline (typically
#        1mmol above), this point is known as"
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("point where we start to accumulate lactate in the legs, this point is known as", `
    $false, $false, $false, $false, $false, $true, 1, $false, `
    "point where lactate levels increase above the baseline (typically 1mmol above), this point is known as", 2) | Out-Null

# Move the (hidden) _GoBack bookmark so it sits between "(typically " and
# "1mmol above)" -- matching where the live edit actually happened.
$rTypically = $d.Content
$rTypically.Find.Execute("(typically 1mmol above)") | Out-Null
$goBackPos = $rTypically.Start + ("(typically ").Length
$goBackRange = $d.Range($goBackPos, $goBackPos)
$d.Bookmarks.Add("_GoBack", $goBackRange) | Out-Null

# ---------------------------------------------------------------------------
# 3. "onset of blood lactate accumulation (OBLA)" -> "aerobic threshold"
#    ("aerobic" gains an underline in addition to the existing bold)
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("onset of blood lactate accumulation (OBLA)", $false, $false, $false, $false, $false, $true, 1, $false, "aerobic threshold", 2) | Out-Null

$rAerobic = $d.Content
$rAerobic.Find.Execute([char]0x201C + "LT1" + [char]0x201D + " or the aerobic threshold") | Out-Null
$aerobicCtx = $rAerobic.Text
$aerobicCtxStart = $rAerobic.Start
$aerobicOffset = $aerobicCtx.IndexOf("aerobic")
$aerobicWordRange = $d.Range($aerobicCtxStart + $aerobicOffset, $aerobicCtxStart + $aerobicOffset + ("aerobic").Length)
$aerobicWordRange.Font.Underline = 1

# ---------------------------------------------------------------------------
# 4. "working, definitely above a" -> "working, but no more than a"
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("working, definitely above a", $false, $false, $false, $false, $false, $true, 1, $false, "working, but no more than a", 2) | Out-Null

# ---------------------------------------------------------------------------
# 5. "the maximal lactate steady state (MLSS) and is closely related to"
#    -> "the anaerobic threshold and more technically the maximal lactate
#         steady state (MLSS). The intensity / power at this point is
#         closely related to"
#    with "anaerobic threshold" bold, "anaerobic" also underlined, and
#    "maximal lactate steady state (MLSS)." bold.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("the maximal lactate steady state (MLSS) and is closely related to", `
    $false, $false, $false, $false, $false, $true, 1, $false, `
    "the anaerobic threshold and more technically the maximal lactate steady state (MLSS). The intensity / power at this point is closely related to", 2) | Out-Null

$rAnaerobic = $d.Content
$rAnaerobic.Find.Execute([char]0x201C + "LT2" + [char]0x201D + " or the anaerobic threshold and more technically the maximal lactate steady state (MLSS).") | Out-Null
$anaerobicCtx = $rAnaerobic.Text
$anaerobicCtxStart = $rAnaerobic.Start

$anaerobicThresholdOffset = $anaerobicCtx.IndexOf("anaerobic threshold")
$anaerobicThresholdRange = $d.Range($anaerobicCtxStart + $anaerobicThresholdOffset, $anaerobicCtxStart + $anaerobicThresholdOffset + ("anaerobic threshold").Length)
$anaerobicThresholdRange.Font.Bold = 1

$anaerobicWordRange = $d.Range($anaerobicCtxStart + $anaerobicThresholdOffset, $anaerobicCtxStart + $anaerobicThresholdOffset + ("anaerobic").Length)
$anaerobicWordRange.Font.Underline = 1

$mlssOffset = $anaerobicCtx.IndexOf("maximal lactate steady state (MLSS).")
$mlssRange = $d.Range($anaerobicCtxStart + $mlssOffset, $anaerobicCtxStart + $mlssOffset + ("maximal lactate steady state (MLSS).").Length)
$mlssRange.Font.Bold = 1

# ---------------------------------------------------------------------------
# 6. " these are the destinations for " + "that " -> merge into a single run
#    (text is unchanged; this is where the _GoBack bookmark used to sit).
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("destinations for that", $false, $false, $false, $false, $false, $true, 1, $false, "destinations for that", 2) | Out-Null

Write-Output "edits applied"
